$wb = $excel.ActiveWorkbook

# --- Sheet "sets": update row 35 (D35: 7 -> 12, E35: 5 -> 10) ---
$wsSets = $wb.Worksheets.Item("sets")
$wsSets.Range("D35").Value = 12
$wsSets.Range("E35").Value = 10

# --- Sheet "rallies": append new rows 257-266 (A257:P266) ---
$wsRallies = $wb.Worksheets.Item("rallies")

# Columns: A rally_id, B match_id, C set_number, D rally_no, E side, F position,
#          G player_number, H action, I result, J who_scored, K score_home,
#          L score_away, M raw_text, N position_zone, O pos_fb, P frente_fundo
# Column F (position) is always blank for these rows.
$data = @(
    @(266, 10, 1, 13, "NOS", 6, "PIPE",      "PONTO", "NOS", 8,  5,  "1 6 pi",   "FRENTE", "FRENTE", "FRENTE"),
    @(267, 10, 1, 14, "NOS", 6, "PIPE",      "ERRO",  "ADV", 8,  6,  "1 6 pi e", "FRENTE", "FRENTE", "FRENTE"),
    @(268, 10, 1, 15, "NOS", 6, "PIPE",      "ERRO",  "ADV", 8,  7,  "1 6 pi e", "FRENTE", "FRENTE", "FRENTE"),
    @(269, 10, 1, 16, "NOS", 6, "PIPE",      "PONTO", "NOS", 9,  7,  "1 6 pi",   "FRENTE", "FRENTE", "FRENTE"),
    @(270, 10, 1, 17, "NOS", 6, "PIPE",      "ERRO",  "ADV", 9,  8,  "1 6 pi e", "FRENTE", "FRENTE", "FRENTE"),
    @(271, 10, 1, 18, "NOS", 7, "RECEPÇÃO",  "PONTO", "NOS", 10, 8,  "1 7 re",   "FRENTE", "FRENTE", "FRENTE"),
    @(272, 10, 1, 19, "NOS", 7, "RECEPÇÃO",  "ERRO",  "ADV", 10, 9,  "1 7 re e", "FRENTE", "FRENTE", "FRENTE"),
    @(273, 10, 1, 20, "NOS", 7, "RECEPÇÃO",  "PONTO", "NOS", 11, 9,  "1 7 re",   "FRENTE", "FRENTE", "FRENTE"),
    @(274, 10, 1, 21, "NOS", 7, "RECEPÇÃO",  "ERRO",  "ADV", 11, 10, "1 7 re e", "FRENTE", "FRENTE", "FRENTE"),
    @(275, 10, 1, 22, "NOS", 7, "RECEPÇÃO",  "PONTO", "NOS", 12, 10, "1 7 re",   "FRENTE", "FRENTE", "FRENTE")
)

$startRow = 257
$lastExistingRow = 256

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    # Copy the last existing data row down into the new row first. This
    # creates all 16 cells (A:P) in the new row - including an empty
    # placeholder for column F (position), matching the blank "position"
    # column used throughout this sheet - before the real values are
    # written on top of it.
    $srcRow = $wsRallies.Range("A" + $lastExistingRow + ":P" + $lastExistingRow)
    $dstRow = $wsRallies.Range("A" + $r + ":P" + $r)
    $srcRow.Copy($dstRow)

    $wsRallies.Cells.Item($r, 1).Value = $row[0]
    $wsRallies.Cells.Item($r, 2).Value = $row[1]
    $wsRallies.Cells.Item($r, 3).Value = $row[2]
    $wsRallies.Cells.Item($r, 4).Value = $row[3]
    $wsRallies.Cells.Item($r, 5).Value = $row[4]
    # column 6 (F) intentionally left as the blank cell created by the copy above
    $wsRallies.Cells.Item($r, 7).Value = $row[5]
    $wsRallies.Cells.Item($r, 8).Value = $row[6]
    $wsRallies.Cells.Item($r, 9).Value = $row[7]
    $wsRallies.Cells.Item($r, 10).Value = $row[8]
    $wsRallies.Cells.Item($r, 11).Value = $row[9]
    $wsRallies.Cells.Item($r, 12).Value = $row[10]
    $wsRallies.Cells.Item($r, 13).Value = $row[11]
    $wsRallies.Cells.Item($r, 14).Value = $row[12]
    $wsRallies.Cells.Item($r, 15).Value = $row[13]
    $wsRallies.Cells.Item($r, 16).Value = $row[14]
}
